$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-order existing match rows (columns F:V) to match canonical fixture order ---
$ws.Range("F2:V2").Copy()
$ws.Range("F100:V100").PasteSpecial(-4163)
$ws.Range("F6:V6").Copy()
$ws.Range("F2:V2").PasteSpecial(-4163)
$ws.Range("F8:V8").Copy()
$ws.Range("F6:V6").PasteSpecial(-4163)
$ws.Range("F9:V9").Copy()
$ws.Range("F8:V8").PasteSpecial(-4163)
$ws.Range("F100:V100").Copy()
$ws.Range("F9:V9").PasteSpecial(-4163)

$ws.Range("F12:V12").Copy()
$ws.Range("F100:V100").PasteSpecial(-4163)
$ws.Range("F18:V18").Copy()
$ws.Range("F12:V12").PasteSpecial(-4163)
$ws.Range("F100:V100").Copy()
$ws.Range("F18:V18").PasteSpecial(-4163)

$ws.Range("F32:V32").Copy()
$ws.Range("F100:V100").PasteSpecial(-4163)
$ws.Range("F33:V33").Copy()
$ws.Range("F32:V32").PasteSpecial(-4163)
$ws.Range("F34:V34").Copy()
$ws.Range("F33:V33").PasteSpecial(-4163)
$ws.Range("F36:V36").Copy()
$ws.Range("F34:V34").PasteSpecial(-4163)
$ws.Range("F39:V39").Copy()
$ws.Range("F36:V36").PasteSpecial(-4163)
$ws.Range("F100:V100").Copy()
$ws.Range("F39:V39").PasteSpecial(-4163)

$ws.Range("F35:V35").Copy()
$ws.Range("F100:V100").PasteSpecial(-4163)
$ws.Range("F37:V37").Copy()
$ws.Range("F35:V35").PasteSpecial(-4163)
$ws.Range("F100:V100").Copy()
$ws.Range("F37:V37").PasteSpecial(-4163)

$ws.Range("F52:V52").Copy()
$ws.Range("F100:V100").PasteSpecial(-4163)
$ws.Range("F58:V58").Copy()
$ws.Range("F52:V52").PasteSpecial(-4163)
$ws.Range("F57:V57").Copy()
$ws.Range("F58:V58").PasteSpecial(-4163)
$ws.Range("F56:V56").Copy()
$ws.Range("F57:V57").PasteSpecial(-4163)
$ws.Range("F59:V59").Copy()
$ws.Range("F56:V56").PasteSpecial(-4163)
$ws.Range("F55:V55").Copy()
$ws.Range("F59:V59").PasteSpecial(-4163)
$ws.Range("F54:V54").Copy()
$ws.Range("F55:V55").PasteSpecial(-4163)
$ws.Range("F53:V53").Copy()
$ws.Range("F54:V54").PasteSpecial(-4163)
$ws.Range("F100:V100").Copy()
$ws.Range("F53:V53").PasteSpecial(-4163)

$ws.Range("F100:V100").ClearContents()

# --- Append new match rows 61-69 (Indice 60-68) ---
$ws.Cells.Item(61, 1).Value = 60
$ws.Cells.Item(61, 2).Value = "france"
$ws.Cells.Item(61, 3).Value = "ligue-2"
$ws.Cells.Item(61, 4).Value = "2023-2024"
$ws.Cells.Item(61, 5).Value = 45192.625
$ws.Cells.Item(61, 6).Value = "Troyes"
$ws.Cells.Item(61, 7).Value = 1
$ws.Cells.Item(61, 8).Value = "Auxerre"
$ws.Cells.Item(61, 9).Value = 2
$ws.Cells.Item(61, 10).Value = 2.98
$ws.Cells.Item(61, 11).Value = "19/09/2023 06:42"
$ws.Cells.Item(61, 12).Value = 4.19
$ws.Cells.Item(61, 13).Value = "23/09/2023 14:59"
$ws.Cells.Item(61, 14).Value = 3.39
$ws.Cells.Item(61, 15).Value = "19/09/2023 06:42"
$ws.Cells.Item(61, 16).Value = 3.92
$ws.Cells.Item(61, 17).Value = "23/09/2023 14:59"
$ws.Cells.Item(61, 18).Value = 2.32
$ws.Cells.Item(61, 19).Value = "19/09/2023 06:42"
$ws.Cells.Item(61, 20).Value = 1.85
$ws.Cells.Item(61, 21).Value = "23/09/2023 14:59"
$ws.Cells.Item(61, 22).Value = "https://www.betexplorer.com/football/france/ligue-2/troyes-auxerre/IwSsVvLj/"
$ws.Range("A60").Copy()
$ws.Range("A61").PasteSpecial(-4122)
$ws.Range("E60").Copy()
$ws.Range("E61").PasteSpecial(-4122)

$ws.Cells.Item(62, 1).Value = 61
$ws.Cells.Item(62, 2).Value = "france"
$ws.Cells.Item(62, 3).Value = "ligue-2"
$ws.Cells.Item(62, 4).Value = "2023-2024"
$ws.Cells.Item(62, 5).Value = 45192.79166666666
$ws.Cells.Item(62, 6).Value = "AC Ajaccio"
$ws.Cells.Item(62, 7).Value = 2
$ws.Cells.Item(62, 8).Value = "Paris FC"
$ws.Cells.Item(62, 9).Value = 1
$ws.Cells.Item(62, 10).Value = 1.97
$ws.Cells.Item(62, 11).Value = "19/09/2023 06:42"
$ws.Cells.Item(62, 12).Value = 2.14
$ws.Cells.Item(62, 13).Value = "23/09/2023 18:58"
$ws.Cells.Item(62, 14).Value = 3.33
$ws.Cells.Item(62, 15).Value = "19/09/2023 06:42"
$ws.Cells.Item(62, 16).Value = 3.16
$ws.Cells.Item(62, 17).Value = "23/09/2023 18:58"
$ws.Cells.Item(62, 18).Value = 3.93
$ws.Cells.Item(62, 19).Value = "19/09/2023 06:42"
$ws.Cells.Item(62, 20).Value = 3.99
$ws.Cells.Item(62, 21).Value = "23/09/2023 18:58"
$ws.Cells.Item(62, 22).Value = "https://www.betexplorer.com/football/france/ligue-2/ac-ajaccio-paris-fc/SfUaw3DA/"
$ws.Range("A60").Copy()
$ws.Range("A62").PasteSpecial(-4122)
$ws.Range("E60").Copy()
$ws.Range("E62").PasteSpecial(-4122)

$ws.Cells.Item(63, 1).Value = 62
$ws.Cells.Item(63, 2).Value = "france"
$ws.Cells.Item(63, 3).Value = "ligue-2"
$ws.Cells.Item(63, 4).Value = "2023-2024"
$ws.Cells.Item(63, 5).Value = 45192.79166666666
$ws.Cells.Item(63, 6).Value = "Amiens"
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = "Valenciennes"
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 1.78
$ws.Cells.Item(63, 11).Value = "16/09/2023 18:12"
$ws.Cells.Item(63, 12).Value = 1.81
$ws.Cells.Item(63, 13).Value = "23/09/2023 18:54"
$ws.Cells.Item(63, 14).Value = 3.7
$ws.Cells.Item(63, 15).Value = "16/09/2023 18:12"
$ws.Cells.Item(63, 16).Value = 3.53
$ws.Cells.Item(63, 17).Value = "23/09/2023 18:54"
$ws.Cells.Item(63, 18).Value = 4.3
$ws.Cells.Item(63, 19).Value = "16/09/2023 18:12"
$ws.Cells.Item(63, 20).Value = 4.98
$ws.Cells.Item(63, 21).Value = "23/09/2023 18:54"
$ws.Cells.Item(63, 22).Value = "https://www.betexplorer.com/football/france/ligue-2/amiens-sc-valenciennes/YXUevqb4/"
$ws.Range("A60").Copy()
$ws.Range("A63").PasteSpecial(-4122)
$ws.Range("E60").Copy()
$ws.Range("E63").PasteSpecial(-4122)

$ws.Cells.Item(64, 1).Value = 63
$ws.Cells.Item(64, 2).Value = "france"
$ws.Cells.Item(64, 3).Value = "ligue-2"
$ws.Cells.Item(64, 4).Value = "2023-2024"
$ws.Cells.Item(64, 5).Value = 45192.79166666666
$ws.Cells.Item(64, 6).Value = "Angers"
$ws.Cells.Item(64, 7).Value = 2
$ws.Cells.Item(64, 8).Value = "Bastia"
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 2.11
$ws.Cells.Item(64, 11).Value = "19/09/2023 06:42"
$ws.Cells.Item(64, 12).Value = 2.02
$ws.Cells.Item(64, 13).Value = "23/09/2023 18:52"
$ws.Cells.Item(64, 14).Value = 3.31
$ws.Cells.Item(64, 15).Value = "19/09/2023 06:42"
$ws.Cells.Item(64, 16).Value = 3.3
$ws.Cells.Item(64, 17).Value = "23/09/2023 18:52"
$ws.Cells.Item(64, 18).Value = 3.73
$ws.Cells.Item(64, 19).Value = "19/09/2023 06:42"
$ws.Cells.Item(64, 20).Value = 4.24
$ws.Cells.Item(64, 21).Value = "23/09/2023 18:52"
$ws.Cells.Item(64, 22).Value = "https://www.betexplorer.com/football/france/ligue-2/angers-bastia/OMN8ysrN/"
$ws.Range("A60").Copy()
$ws.Range("A64").PasteSpecial(-4122)
$ws.Range("E60").Copy()
$ws.Range("E64").PasteSpecial(-4122)

$ws.Cells.Item(65, 1).Value = 64
$ws.Cells.Item(65, 2).Value = "france"
$ws.Cells.Item(65, 3).Value = "ligue-2"
$ws.Cells.Item(65, 4).Value = "2023-2024"
$ws.Cells.Item(65, 5).Value = 45192.79166666666
$ws.Cells.Item(65, 6).Value = "Concarneau"
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 8).Value = "St Etienne"
$ws.Cells.Item(65, 9).Value = 1
$ws.Cells.Item(65, 10).Value = 3.05
$ws.Cells.Item(65, 11).Value = "19/09/2023 06:42"
$ws.Cells.Item(65, 12).Value = 3.36
$ws.Cells.Item(65, 13).Value = "23/09/2023 18:36"
$ws.Cells.Item(65, 14).Value = 3.34
$ws.Cells.Item(65, 15).Value = "19/09/2023 06:42"
$ws.Cells.Item(65, 16).Value = 3.52
$ws.Cells.Item(65, 17).Value = "23/09/2023 18:36"
$ws.Cells.Item(65, 18).Value = 2.31
$ws.Cells.Item(65, 19).Value = "19/09/2023 06:42"
$ws.Cells.Item(65, 20).Value = 2.21
$ws.Cells.Item(65, 21).Value = "23/09/2023 18:36"
$ws.Cells.Item(65, 22).Value = "https://www.betexplorer.com/football/france/ligue-2/concarneau-st-etienne/bZ4UhrDc/"
$ws.Range("A60").Copy()
$ws.Range("A65").PasteSpecial(-4122)
$ws.Range("E60").Copy()
$ws.Range("E65").PasteSpecial(-4122)

$ws.Cells.Item(66, 1).Value = 65
$ws.Cells.Item(66, 2).Value = "france"
$ws.Cells.Item(66, 3).Value = "ligue-2"
$ws.Cells.Item(66, 4).Value = "2023-2024"
$ws.Cells.Item(66, 5).Value = 45192.79166666666
$ws.Cells.Item(66, 6).Value = "Dunkerque"
$ws.Cells.Item(66, 7).Value = 1
$ws.Cells.Item(66, 8).Value = "Rodez"
$ws.Cells.Item(66, 9).Value = 2
$ws.Cells.Item(66, 10).Value = 2.25
$ws.Cells.Item(66, 11).Value = "19/09/2023 06:42"
$ws.Cells.Item(66, 12).Value = 2.44
$ws.Cells.Item(66, 13).Value = "23/09/2023 18:28"
$ws.Cells.Item(66, 14).Value = 3.27
$ws.Cells.Item(66, 15).Value = "19/09/2023 06:42"
$ws.Cells.Item(66, 16).Value = 3.11
$ws.Cells.Item(66, 17).Value = "23/09/2023 18:28"
$ws.Cells.Item(66, 18).Value = 3.22
$ws.Cells.Item(66, 19).Value = "19/09/2023 06:42"
$ws.Cells.Item(66, 20).Value = 3.32
$ws.Cells.Item(66, 21).Value = "23/09/2023 18:28"
$ws.Cells.Item(66, 22).Value = "https://www.betexplorer.com/football/france/ligue-2/dunkerque-rodez/Kr6QgOci/"
$ws.Range("A60").Copy()
$ws.Range("A66").PasteSpecial(-4122)
$ws.Range("E60").Copy()
$ws.Range("E66").PasteSpecial(-4122)

$ws.Cells.Item(67, 1).Value = 66
$ws.Cells.Item(67, 2).Value = "france"
$ws.Cells.Item(67, 3).Value = "ligue-2"
$ws.Cells.Item(67, 4).Value = "2023-2024"
$ws.Cells.Item(67, 5).Value = 45192.79166666666
$ws.Cells.Item(67, 6).Value = "Grenoble"
$ws.Cells.Item(67, 7).Value = 2
$ws.Cells.Item(67, 8).Value = "Quevilly Rouen"
$ws.Cells.Item(67, 9).Value = 1
$ws.Cells.Item(67, 10).Value = 2.76
$ws.Cells.Item(67, 11).Value = "16/09/2023 18:12"
$ws.Cells.Item(67, 12).Value = 2.3
$ws.Cells.Item(67, 13).Value = "23/09/2023 18:37"
$ws.Cells.Item(67, 14).Value = 3.15
$ws.Cells.Item(67, 15).Value = "16/09/2023 18:12"
$ws.Cells.Item(67, 16).Value = 3.14
$ws.Cells.Item(67, 17).Value = "23/09/2023 18:37"
$ws.Cells.Item(67, 18).Value = 2.76
$ws.Cells.Item(67, 19).Value = "16/09/2023 18:12"
$ws.Cells.Item(67, 20).Value = 3.57
$ws.Cells.Item(67, 21).Value = "23/09/2023 18:37"
$ws.Cells.Item(67, 22).Value = "https://www.betexplorer.com/football/france/ligue-2/grenoble-quevilly/S89Yi2S3/"
$ws.Range("A60").Copy()
$ws.Range("A67").PasteSpecial(-4122)
$ws.Range("E60").Copy()
$ws.Range("E67").PasteSpecial(-4122)

$ws.Cells.Item(68, 1).Value = 67
$ws.Cells.Item(68, 2).Value = "france"
$ws.Cells.Item(68, 3).Value = "ligue-2"
$ws.Cells.Item(68, 4).Value = "2023-2024"
$ws.Cells.Item(68, 5).Value = 45192.79166666666
$ws.Cells.Item(68, 6).Value = "Laval"
$ws.Cells.Item(68, 7).Value = 2
$ws.Cells.Item(68, 8).Value = "Guingamp"
$ws.Cells.Item(68, 9).Value = 1
$ws.Cells.Item(68, 10).Value = 2.46
$ws.Cells.Item(68, 11).Value = "18/09/2023 20:13"
$ws.Cells.Item(68, 12).Value = 2.39
$ws.Cells.Item(68, 13).Value = "23/09/2023 18:52"
$ws.Cells.Item(68, 14).Value = 3.13
$ws.Cells.Item(68, 15).Value = "18/09/2023 20:13"
$ws.Cells.Item(68, 16).Value = 3.02
$ws.Cells.Item(68, 17).Value = "23/09/2023 18:55"
$ws.Cells.Item(68, 18).Value = 3
$ws.Cells.Item(68, 19).Value = "18/09/2023 20:13"
$ws.Cells.Item(68, 20).Value = 3.53
$ws.Cells.Item(68, 21).Value = "23/09/2023 18:52"
$ws.Cells.Item(68, 22).Value = "https://www.betexplorer.com/football/france/ligue-2/laval-guingamp/Kh9vWK5p/"
$ws.Range("A60").Copy()
$ws.Range("A68").PasteSpecial(-4122)
$ws.Range("E60").Copy()
$ws.Range("E68").PasteSpecial(-4122)

$ws.Cells.Item(69, 1).Value = 68
$ws.Cells.Item(69, 2).Value = "france"
$ws.Cells.Item(69, 3).Value = "ligue-2"
$ws.Cells.Item(69, 4).Value = "2023-2024"
$ws.Cells.Item(69, 5).Value = 45192.79166666666
$ws.Cells.Item(69, 6).Value = "Pau FC"
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = "Annecy"
$ws.Cells.Item(69, 9).Value = 3
$ws.Cells.Item(69, 10).Value = 2.11
$ws.Cells.Item(69, 11).Value = "16/09/2023 18:12"
$ws.Cells.Item(69, 12).Value = 2.68
$ws.Cells.Item(69, 13).Value = "23/09/2023 18:58"
$ws.Cells.Item(69, 14).Value = 3.38
$ws.Cells.Item(69, 15).Value = "16/09/2023 18:12"
$ws.Cells.Item(69, 16).Value = 3.16
$ws.Cells.Item(69, 17).Value = "23/09/2023 18:54"
$ws.Cells.Item(69, 18).Value = 3.43
$ws.Cells.Item(69, 19).Value = "16/09/2023 18:12"
$ws.Cells.Item(69, 20).Value = 2.91
$ws.Cells.Item(69, 21).Value = "23/09/2023 18:58"
$ws.Cells.Item(69, 22).Value = "https://www.betexplorer.com/football/france/ligue-2/pau-annecy/6oT3xNSG/"
$ws.Range("A60").Copy()
$ws.Range("A69").PasteSpecial(-4122)
$ws.Range("E60").Copy()
$ws.Range("E69").PasteSpecial(-4122)

